$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping_ind")

# Remove "/RME" from the steel description in C2
$cell = $ws.Range("C2")
$cell.Value = "40% S/LFM+CDN/H:1`n20% S+SL/LFM+CDN/H:1`n7% S/LFBR+CDN/H:1`n11% CR/LFM+CDN/H:1`n22% CR+PC/LFM+CDN/H:1"

# Wrap text for the updated cell and grow the row to fit the content
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 256

# Update the active selection, matching the saved view state
$ws.Range("B11").Select()
